$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 214, shifting existing rows 214..296 down to 215..297
$ws.Rows.Item(214).Insert()

# Populate the newly inserted row 214 with the new price record
$ws.Cells.Item(214, 1).Value = 7
$ws.Cells.Item(214, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(214, 3).Value = "Ñuble"
$ws.Cells.Item(214, 4).Value = 44704
$ws.Cells.Item(214, 5).Value = 16
$ws.Cells.Item(214, 6).Value = 100114013
$ws.Cells.Item(214, 7).Value = "Zanahoria"
$ws.Cells.Item(214, 8).Value = "Sin especificar"
$ws.Cells.Item(214, 9).Value = "Primera"
$ws.Cells.Item(214, 10).Value = 120
$ws.Cells.Item(214, 11).Value = 5000
$ws.Cells.Item(214, 12).Value = 5500
$ws.Cells.Item(214, 13).Value = 5250
$ws.Cells.Item(214, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(214, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(214, 16).Value = 262
$ws.Cells.Item(214, 17).Value = 20
$ws.Cells.Item(214, 18).Value = "Hortaliza"
